$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 276, shifting existing rows 276.. down by one.
$ws.Rows("276:276").Insert()

# Populate the newly inserted row 276 with the new data record.
$ws.Cells.Item(276, 1).Value = 6
$ws.Cells.Item(276, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(276, 3).Value = "Metropolitana"
$ws.Cells.Item(276, 4).Value = 44825
$ws.Cells.Item(276, 5).Value = 13
$ws.Cells.Item(276, 6).Value = "Fruta"
$ws.Cells.Item(276, 7).Value = 100101
$ws.Cells.Item(276, 8).Value = "Berries"
$ws.Cells.Item(276, 9).Value = 100101001
$ws.Cells.Item(276, 10).Value = "Arándano (blue)"
$ws.Cells.Item(276, 11).Value = "Sin especificar"
$ws.Cells.Item(276, 12).Value = "Primera"
$ws.Cells.Item(276, 13).Value = 490
$ws.Cells.Item(276, 14).Value = 6000
$ws.Cells.Item(276, 15).Value = 7000
$ws.Cells.Item(276, 16).Value = 6500
$ws.Cells.Item(276, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(276, 18).Value = "Perú"
$ws.Cells.Item(276, 19).Value = 4333
$ws.Cells.Item(276, 20).Value = 1.5

# Apply the same date-style formatting used by column D elsewhere (style index 2).
$ws.Cells.Item(276, 4).NumberFormat = $ws.Cells.Item(277, 4).NumberFormat
